# This workbook (Mobile_B.xlsx) stores a team-specific state-transition
# probability matrix used by a simulation ("simulate game logic").
# The commit "added more games, sped up simulate game logic, and drafted
# optimization logic" re-ran the simulation with a larger number of games,
# which changed the observed transition frequencies (counts / total) in
# several rows of the matrix. Only the probability values themselves
# change; no rows/columns/headers/styles were added, removed, or moved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2068965517241379
$ws.Range("C2").Value = 0.5517241379310345
$ws.Range("P2").Value = 0.1379310344827586
$ws.Range("S2").Value = 0.103448275862069
$ws.Range("B3").Value = 0.05555555555555555
$ws.Range("C3").Value = 0.05555555555555555
$ws.Range("P3").Value = 0.6111111111111112
$ws.Range("S3").Value = 0.2777777777777778
$ws.Range("P4").Value = 0.7142857142857143
$ws.Range("S4").Value = 0.2857142857142857
$ws.Range("B6").Value = 0.09090909090909091
$ws.Range("J6").Value = 0.1818181818181818
$ws.Range("Q6").Value = 0.2727272727272727
$ws.Range("R6").Value = 0.09090909090909091
$ws.Range("S6").Value = 0.3636363636363636
$ws.Range("B7").Value = 0.125
$ws.Range("D7").Value = 0.125
$ws.Range("O7").Value = 0.125
$ws.Range("Q7").Value = 0.125
$ws.Range("B8").Value = 0.06818181818181818
$ws.Range("D8").Value = 0.09090909090909091
$ws.Range("F8").Value = 0.06818181818181818
$ws.Range("J8").Value = 0.1136363636363636
$ws.Range("Q8").Value = 0.09090909090909091
$ws.Range("R8").Value = 0.06818181818181818
$ws.Range("B9").Value = 0.2727272727272727
$ws.Range("J9").Value = 0.09090909090909091
$ws.Range("Q9").Value = 0.1818181818181818
$ws.Range("R9").Value = 0.2727272727272727
$ws.Range("S9").Value = 0.1818181818181818
$ws.Range("B10").Value = 0.1477272727272727
$ws.Range("D10").Value = 0.02272727272727273
$ws.Range("F10").Value = 0.04545454545454546
$ws.Range("J10").Value = 0.09090909090909091
$ws.Range("O10").Value = 0.02272727272727273
$ws.Range("Q10").Value = 0.1931818181818182
$ws.Range("R10").Value = 0.05681818181818182
$ws.Range("S10").Value = 0.4204545454545455
$ws.Range("G11").Value = 0.1481481481481481
$ws.Range("J11").Value = 0.1851851851851852
$ws.Range("K11").Value = 0.3333333333333333
$ws.Range("L11").Value = 0.3333333333333333
$ws.Range("G12").Value = 0.5
$ws.Range("J12").Value = 0.4
$ws.Range("L12").Value = 0.1
$ws.Range("H15").Value = 0.07142857142857142
$ws.Range("I15").Value = 0.1428571428571428
$ws.Range("J15").Value = 0.3571428571428572
$ws.Range("K15").Value = 0.07142857142857142
$ws.Range("O15").Value = 0.1428571428571428
$ws.Range("S15").Value = 0.2142857142857143
$ws.Range("H16").Value = 0.1
$ws.Range("I16").Value = 0.05
$ws.Range("J16").Value = 0.6
$ws.Range("K16").Value = 0.2
$ws.Range("O16").Value = 0.05
$ws.Range("H17").Value = 0.2222222222222222
$ws.Range("I17").Value = 0.03703703703703703
$ws.Range("J17").Value = 0.2592592592592592
$ws.Range("K17").Value = 0.1111111111111111
$ws.Range("M17").Value = 0.03703703703703703
$ws.Range("O17").Value = 0.07407407407407407
$ws.Range("S17").Value = 0.2592592592592592
$ws.Range("H18").Value = 0.2857142857142857
$ws.Range("I18").Value = 0.07142857142857142
$ws.Range("K18").Value = 0.1428571428571428
$ws.Range("M18").Value = 0.07142857142857142
$ws.Range("F19").Value = 0.01063829787234043
$ws.Range("H19").Value = 0.3297872340425532
$ws.Range("I19").Value = 0.06382978723404255
$ws.Range("J19").Value = 0.3617021276595745
$ws.Range("K19").Value = 0.0851063829787234
$ws.Range("M19").Value = 0.02127659574468085
$ws.Range("O19").Value = 0.03191489361702127
$ws.Range("S19").Value = 0.09574468085106383
